$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of (row, col, old, new), using Word's 1-based table indexing.
# Rows 1, 5, 9, 13, 17 contain the data; the others are blank spacer rows.
$updates = @(
    @{ Row = 1;  Col = 1; Old = "10÷5="; New = "61÷3=" },
    @{ Row = 1;  Col = 2; Old = "63÷7="; New = "39÷7=" },
    @{ Row = 1;  Col = 3; Old = "98÷3="; New = "32÷3=" },
    @{ Row = 1;  Col = 4; Old = "14÷4="; New = "55÷8=" },
    @{ Row = 1;  Col = 5; Old = "58÷4="; New = "73÷2=" },

    @{ Row = 5;  Col = 1; Old = "25÷6="; New = "29÷3=" },
    @{ Row = 5;  Col = 2; Old = "49÷6="; New = "40÷3=" },
    @{ Row = 5;  Col = 3; Old = "13÷4="; New = "96÷9=" },
    @{ Row = 5;  Col = 4; Old = "51÷6="; New = "54÷6=" },
    @{ Row = 5;  Col = 5; Old = "74÷3="; New = "37÷8=" },

    @{ Row = 9;  Col = 1; Old = "60÷5="; New = "93÷8=" },
    @{ Row = 9;  Col = 2; Old = "26÷9="; New = "76÷2=" },
    @{ Row = 9;  Col = 3; Old = "56÷4="; New = "28÷2=" },
    @{ Row = 9;  Col = 4; Old = "44÷3="; New = "94÷6=" },
    @{ Row = 9;  Col = 5; Old = "26÷4="; New = "86÷4=" },

    @{ Row = 13; Col = 1; Old = "30÷9="; New = "15÷7=" },
    @{ Row = 13; Col = 2; Old = "93÷9="; New = "38÷8=" },
    @{ Row = 13; Col = 3; Old = "69÷7="; New = "43÷7=" },
    @{ Row = 13; Col = 4; Old = "98÷3="; New = "61÷6=" },
    @{ Row = 13; Col = 5; Old = "36÷4="; New = "29÷2=" },

    @{ Row = 17; Col = 1; Old = "29÷8="; New = "32÷4=" },
    @{ Row = 17; Col = 2; Old = "34÷8="; New = "67÷2=" },
    @{ Row = 17; Col = 3; Old = "78÷7="; New = "26÷8=" },
    @{ Row = 17; Col = 4; Old = "80÷8="; New = "59÷3=" },
    @{ Row = 17; Col = 5; Old = "32÷3="; New = "78÷8=" }
)

foreach ($u in $updates) {
    $cellRange = $t.Cell($u.Row, $u.Col).Range
    # $cellRange.Text includes trailing end-of-cell markers, so compare only
    # the leading substring against the expected original value.
    $current = $cellRange.Text
    $currentValue = $current.Substring(0, $current.Length - 2)
    if ($currentValue -ne $u.Old) {
        Write-Host "Unexpected content at Row=$($u.Row) Col=$($u.Col): found '$currentValue', expected '$($u.Old)'"
    }
    $cellRange.Text = $u.New
}
